$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Columns.Item(4).Insert()

# Copy number/date formatting from the (now-shifted) E column onto the new D column
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Populate column D (new period) and refresh D:K with the restated figures
$row = 7
$ws.Cells.Item($row, 4).Value = 43465
$ws.Cells.Item($row, 5).Value = 43100
$ws.Cells.Item($row, 6).Value = 42735
$ws.Cells.Item($row, 7).Value = 42369
$ws.Cells.Item($row, 8).Value = 42004
$ws.Cells.Item($row, 9).Value = 41639
$ws.Cells.Item($row, 10).Value = 41274
$ws.Cells.Item($row, 11).Value = 40908

$row = 8
$ws.Cells.Item($row, 4).Value = 18565400
$ws.Cells.Item($row, 5).Value = 19953000
$ws.Cells.Item($row, 6).Value = 18591100
$ws.Cells.Item($row, 7).Value = 18779400
$ws.Cells.Item($row, 8).Value = 17762900
$ws.Cells.Item($row, 9).Value = 16392000
$ws.Cells.Item($row, 10).Value = 15483800
$ws.Cells.Item($row, 11).Value = 14755400

$row = 9
$ws.Cells.Item($row, 4).Value = 12781900
$ws.Cells.Item($row, 5).Value = 13200600
$ws.Cells.Item($row, 6).Value = 24252100
$ws.Cells.Item($row, 7).Value = 12797900
$ws.Cells.Item($row, 8).Value = 12157600
$ws.Cells.Item($row, 9).Value = 11075500
$ws.Cells.Item($row, 10).Value = 10321200
$ws.Cells.Item($row, 11).Value = 9881700

$row = 10
$ws.Cells.Item($row, 4).Value = 5783600
$ws.Cells.Item($row, 5).Value = 6752400
$ws.Cells.Item($row, 6).Value = -5661000
$ws.Cells.Item($row, 7).Value = 5981500
$ws.Cells.Item($row, 8).Value = 5605300
$ws.Cells.Item($row, 9).Value = 5316400
$ws.Cells.Item($row, 10).Value = 5162600
$ws.Cells.Item($row, 11).Value = 4873700

$row = 12
$ws.Cells.Item($row, 4).Value = 149500
$ws.Cells.Item($row, 5).Value = 146200
$ws.Cells.Item($row, 6).Value = 163600
$ws.Cells.Item($row, 7).Value = 157400
$ws.Cells.Item($row, 8).Value = 137000
$ws.Cells.Item($row, 9).Value = 141200
$ws.Cells.Item($row, 10).Value = 125200
$ws.Cells.Item($row, 11).Value = 130100

$row = 13
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 14
$ws.Cells.Item($row, 4).Value = -748500
$ws.Cells.Item($row, 5).Value = 195500
$ws.Cells.Item($row, 6).Value = -15200
$ws.Cells.Item($row, 7).Value = "NA"
$ws.Cells.Item($row, 8).Value = "NA"
$ws.Cells.Item($row, 9).Value = "NA"
$ws.Cells.Item($row, 10).Value = "NA"
$ws.Cells.Item($row, 11).Value = 0

$row = 15
$ws.Cells.Item($row, 4).Value = 400
$ws.Cells.Item($row, 5).Value = 500
$ws.Cells.Item($row, 6).Value = 800
$ws.Cells.Item($row, 7).Value = "NA"
$ws.Cells.Item($row, 8).Value = "NA"
$ws.Cells.Item($row, 9).Value = "NA"
$ws.Cells.Item($row, 10).Value = "NA"
$ws.Cells.Item($row, 11).Value = "NA"

$row = 17
$ws.Cells.Item($row, 4).Value = 15157000
$ws.Cells.Item($row, 5).Value = 17302400
$ws.Cells.Item($row, 6).Value = 15888300
$ws.Cells.Item($row, 7).Value = 16168800
$ws.Cells.Item($row, 8).Value = 15233300
$ws.Cells.Item($row, 9).Value = 13860500
$ws.Cells.Item($row, 10).Value = 12994600
$ws.Cells.Item($row, 11).Value = 12319900

$row = 18
$ws.Cells.Item($row, 4).Value = 3408400
$ws.Cells.Item($row, 5).Value = 2650600
$ws.Cells.Item($row, 6).Value = 2702800
$ws.Cells.Item($row, 7).Value = 2610600
$ws.Cells.Item($row, 8).Value = 2529600
$ws.Cells.Item($row, 9).Value = 2531400
$ws.Cells.Item($row, 10).Value = 2489200
$ws.Cells.Item($row, 11).Value = 2435500

$row = 20
$ws.Cells.Item($row, 4).Value = 165400
$ws.Cells.Item($row, 5).Value = 57600
$ws.Cells.Item($row, 6).Value = 578500
$ws.Cells.Item($row, 7).Value = 130800
$ws.Cells.Item($row, 8).Value = 94500
$ws.Cells.Item($row, 9).Value = 43700
$ws.Cells.Item($row, 10).Value = 206500
$ws.Cells.Item($row, 11).Value = 70200

$row = 21
$ws.Cells.Item($row, 4).Value = 4389000
$ws.Cells.Item($row, 5).Value = 3535500
$ws.Cells.Item($row, 6).Value = 4153900
$ws.Cells.Item($row, 7).Value = 3548100
$ws.Cells.Item($row, 8).Value = 3410600
$ws.Cells.Item($row, 9).Value = 3304200
$ws.Cells.Item($row, 10).Value = 3373800
$ws.Cells.Item($row, 11).Value = 3160300

$row = 22
$ws.Cells.Item($row, 4).Value = 503200
$ws.Cells.Item($row, 5).Value = 467000
$ws.Cells.Item($row, 6).Value = 986200
$ws.Cells.Item($row, 7).Value = 570000
$ws.Cells.Item($row, 8).Value = 555800
$ws.Cells.Item($row, 9).Value = 502100
$ws.Cells.Item($row, 10).Value = 527900
$ws.Cells.Item($row, 11).Value = 418300

$row = 23
$ws.Cells.Item($row, 4).Value = 3070600
$ws.Cells.Item($row, 5).Value = 2241300
$ws.Cells.Item($row, 6).Value = 2295000
$ws.Cells.Item($row, 7).Value = 2171400
$ws.Cells.Item($row, 8).Value = 2068300
$ws.Cells.Item($row, 9).Value = 2073000
$ws.Cells.Item($row, 10).Value = 2167800
$ws.Cells.Item($row, 11).Value = 2087500

$row = 24
$ws.Cells.Item($row, 4).Value = 573400
$ws.Cells.Item($row, 5).Value = 497100
$ws.Cells.Item($row, 6).Value = 701700
$ws.Cells.Item($row, 7).Value = 698000
$ws.Cells.Item($row, 8).Value = 654800
$ws.Cells.Item($row, 9).Value = 664200
$ws.Cells.Item($row, 10).Value = 679000
$ws.Cells.Item($row, 11).Value = 705600

$row = 25
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 26
$ws.Cells.Item($row, 4).Value = 2497200
$ws.Cells.Item($row, 5).Value = 1744200
$ws.Cells.Item($row, 6).Value = 1593300
$ws.Cells.Item($row, 7).Value = 1473300
$ws.Cells.Item($row, 8).Value = 1413500
$ws.Cells.Item($row, 9).Value = 1408800
$ws.Cells.Item($row, 10).Value = 1488900
$ws.Cells.Item($row, 11).Value = 1381900

$row = 27
$ws.Cells.Item($row, 4).Value = 2223700
$ws.Cells.Item($row, 5).Value = 1435900
$ws.Cells.Item($row, 6).Value = 1283500
$ws.Cells.Item($row, 7).Value = 1155000
$ws.Cells.Item($row, 8).Value = 1172800
$ws.Cells.Item($row, 9).Value = 1245300
$ws.Cells.Item($row, 10).Value = 1331500
$ws.Cells.Item($row, 11).Value = 1257200

$row = 28
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 29
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 30
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 31
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 32
$ws.Cells.Item($row, 4).Value = -165400
$ws.Cells.Item($row, 5).Value = -57600
$ws.Cells.Item($row, 6).Value = -578500
$ws.Cells.Item($row, 7).Value = -130800
$ws.Cells.Item($row, 8).Value = -94500
$ws.Cells.Item($row, 9).Value = -43700
$ws.Cells.Item($row, 10).Value = -206500
$ws.Cells.Item($row, 11).Value = -70200

$row = 33
$ws.Cells.Item($row, 4).Value = 2223700
$ws.Cells.Item($row, 5).Value = 1435900
$ws.Cells.Item($row, 6).Value = 1283500
$ws.Cells.Item($row, 7).Value = 1155000
$ws.Cells.Item($row, 8).Value = 1172800
$ws.Cells.Item($row, 9).Value = 1245300
$ws.Cells.Item($row, 10).Value = 1331500
$ws.Cells.Item($row, 11).Value = 1257200

$row = 34
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 35
$ws.Cells.Item($row, 4).Value = 2223700
$ws.Cells.Item($row, 5).Value = 1435900
$ws.Cells.Item($row, 6).Value = 1283500
$ws.Cells.Item($row, 7).Value = 1155000
$ws.Cells.Item($row, 8).Value = 1172800
$ws.Cells.Item($row, 9).Value = 1245300
$ws.Cells.Item($row, 10).Value = 1331500
$ws.Cells.Item($row, 11).Value = 1257200

$row = 38
$ws.Cells.Item($row, 4).Value = 43465
$ws.Cells.Item($row, 5).Value = 43100
$ws.Cells.Item($row, 6).Value = 42735
$ws.Cells.Item($row, 7).Value = 42369
$ws.Cells.Item($row, 8).Value = 42004
$ws.Cells.Item($row, 9).Value = 41639
$ws.Cells.Item($row, 10).Value = 41274
$ws.Cells.Item($row, 11).Value = 40908

$row = 41
$ws.Cells.Item($row, 4).Value = 2407400
$ws.Cells.Item($row, 5).Value = 1097400
$ws.Cells.Item($row, 6).Value = 838400
$ws.Cells.Item($row, 7).Value = 616500
$ws.Cells.Item($row, 8).Value = 711200
$ws.Cells.Item($row, 9).Value = 766100
$ws.Cells.Item($row, 10).Value = 772000
$ws.Cells.Item($row, 11).Value = 536800

$row = 42
$ws.Cells.Item($row, 4).Value = 134100
$ws.Cells.Item($row, 5).Value = 21300
$ws.Cells.Item($row, 6).Value = 314400
$ws.Cells.Item($row, 7).Value = 322300
$ws.Cells.Item($row, 8).Value = 210400
$ws.Cells.Item($row, 9).Value = 21600
$ws.Cells.Item($row, 10).Value = 23500
$ws.Cells.Item($row, 11).Value = 19400

$row = 43
$ws.Cells.Item($row, 4).Value = 4224900
$ws.Cells.Item($row, 5).Value = 4147000
$ws.Cells.Item($row, 6).Value = 4473900
$ws.Cells.Item($row, 7).Value = 4270600
$ws.Cells.Item($row, 8).Value = 4241300
$ws.Cells.Item($row, 9).Value = 3808500
$ws.Cells.Item($row, 10).Value = 3762000
$ws.Cells.Item($row, 11).Value = 3671900

$row = 44
$ws.Cells.Item($row, 4).Value = 1645700
$ws.Cells.Item($row, 5).Value = 1448200
$ws.Cells.Item($row, 6).Value = 1581800
$ws.Cells.Item($row, 7).Value = 1504300
$ws.Cells.Item($row, 8).Value = 1251600
$ws.Cells.Item($row, 9).Value = 1230900
$ws.Cells.Item($row, 10).Value = 1163300
$ws.Cells.Item($row, 11).Value = 1135700

$row = 45
$ws.Cells.Item($row, 4).Value = 392000
$ws.Cells.Item($row, 5).Value = 1160500
$ws.Cells.Item($row, 6).Value = 997600
$ws.Cells.Item($row, 7).Value = 2665300
$ws.Cells.Item($row, 8).Value = 1130500
$ws.Cells.Item($row, 9).Value = 1226500
$ws.Cells.Item($row, 10).Value = 2206400
$ws.Cells.Item($row, 11).Value = 1321100

$row = 46
$ws.Cells.Item($row, 4).Value = 8804100
$ws.Cells.Item($row, 5).Value = 7151900
$ws.Cells.Item($row, 6).Value = 8206200
$ws.Cells.Item($row, 7).Value = 7594100
$ws.Cells.Item($row, 8).Value = 7545100
$ws.Cells.Item($row, 9).Value = 7053600
$ws.Cells.Item($row, 10).Value = 6874900
$ws.Cells.Item($row, 11).Value = 6684900

$row = 47
$ws.Cells.Item($row, 4).Value = 729000
$ws.Cells.Item($row, 5).Value = 725900
$ws.Cells.Item($row, 6).Value = 762100
$ws.Cells.Item($row, 7).Value = 723400
$ws.Cells.Item($row, 8).Value = 759400
$ws.Cells.Item($row, 9).Value = 745500
$ws.Cells.Item($row, 10).Value = 715100
$ws.Cells.Item($row, 11).Value = 812300

$row = 48
$ws.Cells.Item($row, 4).Value = 4304000
$ws.Cells.Item($row, 5).Value = 7835500
$ws.Cells.Item($row, 6).Value = 4233500
$ws.Cells.Item($row, 7).Value = 7686900
$ws.Cells.Item($row, 8).Value = 3691500
$ws.Cells.Item($row, 9).Value = 3469100
$ws.Cells.Item($row, 10).Value = 6598700
$ws.Cells.Item($row, 11).Value = 3086800

$row = 49
$ws.Cells.Item($row, 4).Value = 14463500
$ws.Cells.Item($row, 5).Value = 15113200
$ws.Cells.Item($row, 6).Value = 16284200
$ws.Cells.Item($row, 7).Value = 16486200
$ws.Cells.Item($row, 8).Value = 15653500
$ws.Cells.Item($row, 9).Value = 13930700
$ws.Cells.Item($row, 10).Value = 14408700
$ws.Cells.Item($row, 11).Value = 11589400

$row = 50
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 51
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 52
$ws.Cells.Item($row, 4).Value = 1142900
$ws.Cells.Item($row, 5).Value = 813600
$ws.Cells.Item($row, 6).Value = 733700
$ws.Cells.Item($row, 7).Value = 744200
$ws.Cells.Item($row, 8).Value = 901900
$ws.Cells.Item($row, 9).Value = 741300
$ws.Cells.Item($row, 10).Value = 548200
$ws.Cells.Item($row, 11).Value = 754500

$row = 53
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 54
$ws.Cells.Item($row, 4).Value = 29443600
$ws.Cells.Item($row, 5).Value = 26956100
$ws.Cells.Item($row, 6).Value = 30219600
$ws.Cells.Item($row, 7).Value = 28459600
$ws.Cells.Item($row, 8).Value = 28551400
$ws.Cells.Item($row, 9).Value = 25940300
$ws.Cells.Item($row, 10).Value = 25049500
$ws.Cells.Item($row, 11).Value = 22927900

$row = 57
$ws.Cells.Item($row, 4).Value = 719500
$ws.Cells.Item($row, 5).Value = 662500
$ws.Cells.Item($row, 6).Value = 680700
$ws.Cells.Item($row, 7).Value = 704400
$ws.Cells.Item($row, 8).Value = 643100
$ws.Cells.Item($row, 9).Value = 608800
$ws.Cells.Item($row, 10).Value = 698200
$ws.Cells.Item($row, 11).Value = 635500

$row = 58
$ws.Cells.Item($row, 4).Value = 2960800
$ws.Cells.Item($row, 5).Value = 2970500
$ws.Cells.Item($row, 6).Value = 1536100
$ws.Cells.Item($row, 7).Value = 1011900
$ws.Cells.Item($row, 8).Value = 506800
$ws.Cells.Item($row, 9).Value = 752100
$ws.Cells.Item($row, 10).Value = 512300
$ws.Cells.Item($row, 11).Value = 2015000

$row = 59
$ws.Cells.Item($row, 4).Value = 3352600
$ws.Cells.Item($row, 5).Value = 6495400
$ws.Cells.Item($row, 6).Value = 3435000
$ws.Cells.Item($row, 7).Value = 5898200
$ws.Cells.Item($row, 8).Value = 2751600
$ws.Cells.Item($row, 9).Value = 2626600
$ws.Cells.Item($row, 10).Value = 3373100
$ws.Cells.Item($row, 11).Value = 2353000

$row = 60
$ws.Cells.Item($row, 4).Value = 7032900
$ws.Cells.Item($row, 5).Value = 5946400
$ws.Cells.Item($row, 6).Value = 5651800
$ws.Cells.Item($row, 7).Value = 4655600
$ws.Cells.Item($row, 8).Value = 3901500
$ws.Cells.Item($row, 9).Value = 3987500
$ws.Cells.Item($row, 10).Value = 3556700
$ws.Cells.Item($row, 11).Value = 5003500

$row = 61
$ws.Cells.Item($row, 4).Value = 5661000
$ws.Cells.Item($row, 5).Value = 6501800
$ws.Cells.Item($row, 6).Value = 8081200
$ws.Cells.Item($row, 7).Value = 8811500
$ws.Cells.Item($row, 8).Value = 10188000
$ws.Cells.Item($row, 9).Value = 8692000
$ws.Cells.Item($row, 10).Value = 8798500
$ws.Cells.Item($row, 11).Value = 6449900

$row = 62
$ws.Cells.Item($row, 4).Value = 2273700
$ws.Cells.Item($row, 5).Value = 3597500
$ws.Cells.Item($row, 6).Value = 2239400
$ws.Cells.Item($row, 7).Value = 2062200
$ws.Cells.Item($row, 8).Value = 2284900
$ws.Cells.Item($row, 9).Value = 1891400
$ws.Cells.Item($row, 10).Value = 1776800
$ws.Cells.Item($row, 11).Value = 1530600

$row = 63
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 64
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 65
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 66
$ws.Cells.Item($row, 4).Value = 16250700
$ws.Cells.Item($row, 5).Value = 15938000
$ws.Cells.Item($row, 6).Value = 18092300
$ws.Cells.Item($row, 7).Value = 17366300
$ws.Cells.Item($row, 8).Value = 17956100
$ws.Cells.Item($row, 9).Value = 15579200
$ws.Cells.Item($row, 10).Value = 15016100
$ws.Cells.Item($row, 11).Value = 13652900

$row = 68
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 69
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 70
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 5000
$ws.Cells.Item($row, 11).Value = 5200

$row = 71
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 72
$ws.Cells.Item($row, 4).Value = 9909300
$ws.Cells.Item($row, 5).Value = 8007900
$ws.Cells.Item($row, 6).Value = 9915100
$ws.Cells.Item($row, 7).Value = 8831200
$ws.Cells.Item($row, 8).Value = 7971500
$ws.Cells.Item($row, 9).Value = 7155400
$ws.Cells.Item($row, 10).Value = 6242400
$ws.Cells.Item($row, 11).Value = 5456600

$row = 73
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 74
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 75
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 76
$ws.Cells.Item($row, 4).Value = 13192800
$ws.Cells.Item($row, 5).Value = 11018100
$ws.Cells.Item($row, 6).Value = 12127300
$ws.Cells.Item($row, 7).Value = 11093300
$ws.Cells.Item($row, 8).Value = 10595300
$ws.Cells.Item($row, 9).Value = 10361100
$ws.Cells.Item($row, 10).Value = 10028400
$ws.Cells.Item($row, 11).Value = 9269700

$row = 77
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 80
$ws.Cells.Item($row, 4).Value = 43465
$ws.Cells.Item($row, 5).Value = 43100
$ws.Cells.Item($row, 6).Value = 42735
$ws.Cells.Item($row, 7).Value = 42369
$ws.Cells.Item($row, 8).Value = 42004
$ws.Cells.Item($row, 9).Value = 41639
$ws.Cells.Item($row, 10).Value = 41274
$ws.Cells.Item($row, 11).Value = 40908

$row = 81
$ws.Cells.Item($row, 4).Value = 2223700
$ws.Cells.Item($row, 5).Value = 1435900
$ws.Cells.Item($row, 6).Value = 1283500
$ws.Cells.Item($row, 7).Value = 1155000
$ws.Cells.Item($row, 8).Value = 1172800
$ws.Cells.Item($row, 9).Value = 1245300
$ws.Cells.Item($row, 10).Value = 1331500
$ws.Cells.Item($row, 11).Value = 1257200

$row = 83
$ws.Cells.Item($row, 4).Value = 813300
$ws.Cells.Item($row, 5).Value = 825200
$ws.Cells.Item($row, 6).Value = 870600
$ws.Cells.Item($row, 7).Value = 804800
$ws.Cells.Item($row, 8).Value = 784600
$ws.Cells.Item($row, 9).Value = 727300
$ws.Cells.Item($row, 10).Value = 676400
$ws.Cells.Item($row, 11).Value = 654100

$row = 84
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 85
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 86
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 87
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 88
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 89
$ws.Cells.Item($row, 4).Value = 2313400
$ws.Cells.Item($row, 5).Value = 2459200
$ws.Cells.Item($row, 6).Value = 2400900
$ws.Cells.Item($row, 7).Value = 2199200
$ws.Cells.Item($row, 8).Value = 2088500
$ws.Cells.Item($row, 9).Value = 2283000
$ws.Cells.Item($row, 10).Value = 2287800
$ws.Cells.Item($row, 11).Value = 1697900

$row = 91
$ws.Cells.Item($row, 4).Value = -1186300
$ws.Cells.Item($row, 5).Value = -1059700
$ws.Cells.Item($row, 6).Value = -1155600
$ws.Cells.Item($row, 7).Value = -1069200
$ws.Cells.Item($row, 8).Value = -1045300
$ws.Cells.Item($row, 9).Value = -839200
$ws.Cells.Item($row, 10).Value = -757700
$ws.Cells.Item($row, 11).Value = -2797400

$row = 92
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 93
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 94
$ws.Cells.Item($row, 4).Value = -274900
$ws.Cells.Item($row, 5).Value = -1112500
$ws.Cells.Item($row, 6).Value = -1547600
$ws.Cells.Item($row, 7).Value = -1122800
$ws.Cells.Item($row, 8).Value = -3019000
$ws.Cells.Item($row, 9).Value = -1352600
$ws.Cells.Item($row, 10).Value = -2559500
$ws.Cells.Item($row, 11).Value = -2753600

$row = 96
$ws.Cells.Item($row, 4).Value = -364500
$ws.Cells.Item($row, 5).Value = -329800
$ws.Cells.Item($row, 6).Value = -311000
$ws.Cells.Item($row, 7).Value = -295400
$ws.Cells.Item($row, 8).Value = -356700
$ws.Cells.Item($row, 9).Value = -332300
$ws.Cells.Item($row, 10).Value = -304900
$ws.Cells.Item($row, 11).Value = -329400

$row = 97
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 98
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 99
$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0

$row = 100
$ws.Cells.Item($row, 4).Value = -764900
$ws.Cells.Item($row, 5).Value = -896100
$ws.Cells.Item($row, 6).Value = -656100
$ws.Cells.Item($row, 7).Value = -1130400
$ws.Cells.Item($row, 8).Value = 903200
$ws.Cells.Item($row, 9).Value = -906600
$ws.Cells.Item($row, 10).Value = 525500
$ws.Cells.Item($row, 11).Value = 931000

$row = 101
$ws.Cells.Item($row, 4).Value = 36300
$ws.Cells.Item($row, 5).Value = -148600
$ws.Cells.Item($row, 6).Value = 24600
$ws.Cells.Item($row, 7).Value = -40600
$ws.Cells.Item($row, 8).Value = -27600
$ws.Cells.Item($row, 9).Value = -29700
$ws.Cells.Item($row, 10).Value = 5100
$ws.Cells.Item($row, 11).Value = 47700

$row = 102
$ws.Cells.Item($row, 4).Value = 1309900
$ws.Cells.Item($row, 5).Value = 302100
$ws.Cells.Item($row, 6).Value = 221900
$ws.Cells.Item($row, 7).Value = -94600
$ws.Cells.Item($row, 8).Value = -54900
$ws.Cells.Item($row, 9).Value = -5900
$ws.Cells.Item($row, 10).Value = 258900
$ws.Cells.Item($row, 11).Value = -77000
